$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 15.94901084899902
$ws.Range("D2").Value = 123

$ws.Range("C3").Value = 16.72792434692383
$ws.Range("D3").Value = 175

$ws.Range("C4").Value = 16.00790023803711
$ws.Range("D4").Value = 175

$ws.Range("C5").Value = 15.96188545227051
$ws.Range("D5").Value = 123

$ws.Range("C6").Value = 16.29495620727539
$ws.Range("D6").Value = 176
